# Tracking Stunden dem Terminplan angepasst
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

function Set-RowValues($ws, $rangeAddress, $values) {
    $n = $values.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range($rangeAddress).Value = $arr
}

# --- Row 14: "1_Projektmanagement" planned hours per week (Std / W(T)) ---
Set-RowValues $ws "E14:Y14" @(1,2,4,4,3,4,2,2,2,2,3,1,1,1,1,1,1,1,2,3,3)

# --- Row 28: "2_Analyse&Entwurf" planned hours per week (Std / W(T)) ---
Set-RowValues $ws "E28:L28" @(10,15,20,25,25,15,15,10)

# --- Row 62: "4_Realisierung_SW" planned hours per week (Std / W(T)) ---
Set-RowValues $ws "G62:U62" @(5,5,6,6,12,12,15,24,15,24,24,24,24,24,5)

# Update the selected cell on the sheet to match the saved view state
$ws.Range("V62").Select()

$wb.Save()
